$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header for the triple_double_avg column, matching the
# formatting of the preceding header cell (bold font, border, centered)
$ws.Range("J1").Value = "triple_double_avg"
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Populate the new column with "No" for each data row (rows 2-6)
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 10).Value = "No"
}
